# Insert a new weekly price-report row for "Terminal La Palmera de La Serena -
# Cilantro" at row 161 (pushing the former rows 161-204 down to 162-205), and
# populate the new row with its own data. This mirrors the upstream diff,
# which shows row 161's old contents re-appearing (unchanged) at row 162, row
# 162's old contents at row 163, and so on through the former last row (204)
# which now lives at 205 - i.e. a single row insertion, not per-cell edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 161:204 down to 162:205, leaving a blank row 161 behind.
$ws.Rows("161:161").Insert()

# Populate the newly inserted row 161 with the new weekly record. All the
# "constant" columns (A,B,C,E,F,G,H,I,N,O,Q,R) carry the same values as every
# other row in this sheet for this market/product.
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 44964
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100112040
$ws.Range("G161").Value = "Cilantro"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 2000
$ws.Range("K161").Value = 2800
$ws.Range("L161").Value = 3000
$ws.Range("M161").Value = 2900
$ws.Range("N161").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O161").Value = "Provincia del Elquí"
$ws.Range("P161").Value = 1933
$ws.Range("Q161").Value = 1.5
$ws.Range("R161").Value = "Hortaliza"
